# Fruta / hortaliza, semanal
# Insert 5 new weekly price rows for Femacal de La Calera - Cereza,
# shifting the existing data (old rows 716:809) down to (721:814).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows before the current row 716; this pushes the old
# rows 716-809 down to 721-814 and keeps their formatting/values intact.
$ws.Rows("716:720").Insert()

$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 45265, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Brooks", "Especial", 56, 1400, 1400, 1400, "$/kilo (en caja de 15 kilos)", "Región de O'Higgins", 1400, 1),
    @(3, "Femacal de La Calera", "Coquimbo", 45265, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Brooks", "Primera", 50, 1000, 1000, 1000, "$/kilo (en caja de 15 kilos)", "Región de O'Higgins", 1000, 1),
    @(3, "Femacal de La Calera", "Coquimbo", 45265, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Royal Dawn", "Especial", 56, 15000, 15000, 15000, "$/bandeja 10 kilos", "Región de O'Higgins", 1500, 10),
    @(3, "Femacal de La Calera", "Coquimbo", 45265, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Royal Dawn", "Primera", 50, 12000, 12000, 12000, "$/bandeja 10 kilos", "Región de O'Higgins", 1200, 10),
    @(3, "Femacal de La Calera", "Coquimbo", 45265, 5, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Royal Dawn", "Segunda", 48, 10000, 10000, 10000, "$/bandeja 10 kilos", "Región de O'Higgins", 1000, 10)
)

$startRow = 716
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
